# Add a new day's pair of columns (AH = "05-07_A", AI = "05-07_0") to the
# Season_Attack sheet, mirroring the previous day's columns (AF = "05-06_A",
# AG = "05-06_0") the way Excel would when duplicating the last tracked day
# as a starting point for the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 120

# Column indices: AF=32, AG=33 (source) -> AH=34, AI=35 (new)
$colAF = 32
$colAG = 33
$colAH = 34
$colAI = 35

# 1) Copy the whole AF:AG block (values + styles) into AH:AI in one shot.
#    This gives AH the exact same per-row style/value as AF, and AI the
#    exact same per-row style/value as AG (including the header row text,
#    which we fix up right after).
$ws.Range("AF1:AG" + $lastRow).Copy($ws.Range("AH1"))

# 2) Fix up the new header row text (next day's labels).
$ws.Cells.Item(1, $colAH).Value = "05-07_A"
$ws.Cells.Item(1, $colAI).Value = "05-07_0"

# 3) The original "_0" column (AG) used to store its running total as text
#    (inline string). Re-enter each non-blank value through .Value so the
#    engine stores it as a genuine number, matching the new AG semantics -
#    while AI (already populated by the copy above) keeps the old text
#    representation untouched.
for ($r = 2; $r -le $lastRow; $r++) {
    $agCell = $ws.Cells.Item($r, $colAG)
    $current = $agCell.Value2
    if (-not [string]::IsNullOrEmpty($current)) {
        $agCell.Value = $current
    }
}
